# Updated symbol list on Thu Feb 16 09:50:21 UTC 2023 with GitHub Actions
# Refresh Price (col D) and Volume(1h) (col E) figures for the crypto table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. A leading apostrophe forces Excel to keep
# the numeric-looking / percent-looking text as a literal string (matching
# the existing column formatting), and resetting the style back to "Normal"
# avoids leaving a stray quote-prefix style on the cell.
$updates = @{
    "D2" = "321.91"
    "E2" = "8.03%"
    "D3" = "48.26"
    "E3" = "14.32%"
    "D4" = "5.284"
    "E4" = "5.60%"
    "D5" = "0.08104"
    "E5" = "7.77%"
    "D6" = "4.570"
    "E6" = "4.63%"
    "E7" = "2.83%"
    "E8" = "29.07%"
    "E9" = "9.28%"
    "D10" = "0.1946"
    "E10" = "6.05%"
    "D11" = "0.09473"
    "E11" = "6.57%"
    "D12" = "0.04640"
    "E12" = "12.26%"
    "E13" = "-0.09%"
    "D14" = "0.001338"
    "E14" = "3.56%"
    "D15" = "0.04130"
    "E15" = "1.15%"
    "D16" = "0.005801"
    "E16" = "0.03%"
    "D17" = "3.342"
    "E18" = "1.15%"
    "E19" = "2.07%"
    "D20" = "8.084"
    "E20" = "-3.06%"
    "E21" = "3.61%"
    "D23" = "0.001306"
    "E23" = "3.15%"
    "D24" = "0.004251"
    "E24" = "9.23%"
    "E25" = "3.77%"
    "D26" = "0.0003540"
    "E26" = "-4.95%"
    "D38" = "0.02730"
    "E38" = "14.00%"
    "D39" = "0.05892"
    "E39" = "12.82%"
    "E40" = "-6.56%"
    "D41" = "0.007693"
    "E41" = "-0.99%"
    "D42" = "0.1441"
    "E42" = "8.75%"
    "D43" = "0.007717"
    "E43" = "4.17%"
    "E44" = "13.74%"
    "E45" = "6.94%"
    "D46" = "0.00007018"
    "E46" = "6.85%"
    "E47" = "-0.06%"
    "D48" = "0.05377"
    "E48" = "18.98%"
    "D49" = "0.004000"
    "E49" = "-4.82%"
    "E50" = "-0.06%"
    "E51" = "-0.06%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
